$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 243
$ws.Range("C2").Value = 84.08

$ws.Range("C3").Value = 6.23

$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 3.46

$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 2.08

$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 1.73

$ws.Range("B7").Value = 4
$ws.Range("C7").Value = 1.38

$ws.Range("C8").Value = 0.6899999999999999

$ws.Range("C9").Value = 0.35
